$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ranges = @("D31:D32","E31:E32","F31:F32","A31:A32","B31:B32","C31:C32")
foreach ($r in $ranges) {
    $ws.Range($r).Merge()
}
